$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" (last-changed) date for every data row
# (rows 2-238). The whole column was bumped by one day (2023-09-12 ->
# 2023-09-13, serial 45181 -> 45182) as part of an automatic daily update.
$ws.Range("C2:C238").Value = 45182
